$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quarterly_index")

$ws.Range("I1").Value = "struct_pub_share"
$ws.Range("J1").Value = "struct_comp_share"
$ws.Range("K1").Value = "struct_outperf"
$ws.Range("L1").Value = "struct_outperf_yoy"

$ws.Range("I2").Value = 0.108108108
$ws.Range("J2").Value = 0.4
$ws.Range("K2").Value = -0.291891892

$ws.Range("I3").Value = 0.07934579645714285
$ws.Range("J3").Value = 0.4101794638142858
$ws.Range("K3").Value = -0.3308336673571429

$ws.Range("I4").Value = 0.07578112153333333
$ws.Range("J4").Value = 0.3907256863222222
$ws.Range("K4").Value = -0.3149445647888889

$ws.Range("I5").Value = 0.1493226793183093
$ws.Range("J5").Value = 0.3831075597727699
$ws.Range("K5").Value = -0.2337848804544606

$ws.Range("I6").Value = 0.102908277
$ws.Range("J6").Value = 0.402684564
$ws.Range("K6").Value = -0.299776287
$ws.Range("L6").Value = 0.02701135323073678

$ws.Range("I7").Value = 0.05835321335067505
$ws.Range("J7").Value = 0.3872388234499205
$ws.Range("K7").Value = -0.3288856100992454
$ws.Range("L7").Value = -0.005888328335684534

$ws.Range("I8").Value = 0.0693636464975713
$ws.Range("J8").Value = 0.3661496268939584
$ws.Range("K8").Value = -0.296785980396387
$ws.Range("L8").Value = -0.05765644631674705

$ws.Range("I9").Value = 0.07619810651255549
$ws.Range("J9").Value = 0.3642757533389226
$ws.Range("K9").Value = -0.2880776468263671
$ws.Range("L9").Value = 0.2322338650231157

$ws.Range("I10").Value = 0.08095832897003177
$ws.Range("J10").Value = 0.3526428828493964
$ws.Range("K10").Value = -0.2716845538793646
$ws.Range("L10").Value = -0.09370899013315015

$ws.Range("I11").Value = 0.07953337952278711
$ws.Range("J11").Value = 0.3566932889943775
$ws.Range("K11").Value = -0.2771599094715904
$ws.Range("L11").Value = -0.1572756576733355

$ws.Range("I12").Value = 0.0785765924326669
$ws.Range("J12").Value = 0.3501030552297806
$ws.Range("K12").Value = -0.2715264627971137
$ws.Range("L12").Value = -0.08511021162636023

$ws.Range("I13").Value = 0.07997296529438562
$ws.Range("J13").Value = 0.3663393478835924
$ws.Range("K13").Value = -0.2863663825892068
$ws.Range("L13").Value = -0.005940288168876151

$ws.Range("I14").Value = 0.07503083144407637
$ws.Range("J14").Value = 0.3509653241224536
$ws.Range("K14").Value = -0.2759344926783773
$ws.Range("L14").Value = 0.01564291653068994

$ws.Range("I15").Value = 0.0616435026809537
$ws.Range("J15").Value = 0.3031180056651374
$ws.Range("K15").Value = -0.2414745029841837
$ws.Range("L15").Value = -0.1287538538868824

$ws.Range("I16").Value = 0.05387854874033485
$ws.Range("J16").Value = 0.2907285623440602
$ws.Range("K16").Value = -0.2368500136037253
$ws.Range("L16").Value = -0.1277092804737004

$ws.Range("I17").Value = 0.0713397436418451
$ws.Range("J17").Value = 0.3228187098619124
$ws.Range("K17").Value = -0.2514789662200673
$ws.Range("L17").Value = -0.1218279047062082

$ws.Range("I18").Value = 0.07304880225158991
$ws.Range("J18").Value = 0.3089007092112859
$ws.Range("K18").Value = -0.235851906959696
$ws.Range("L18").Value = -0.1452612369320589

$ws.Range("I19").Value = 0.08378966317958914
$ws.Range("J19").Value = 0.2941511762854176
$ws.Range("K19").Value = -0.2103615131058285
$ws.Range("L19").Value = -0.1288458594752468

$ws.Range("I20").Value = 0.09503764666475309
$ws.Range("J20").Value = 0.321860205835205
$ws.Range("K20").Value = -0.2268225591704519
$ws.Range("L20").Value = -0.04233672728450999

$ws.Range("I21").Value = 0.07740772446012015
$ws.Range("J21").Value = 0.2912632172066835
$ws.Range("K21").Value = -0.2138554927465634
$ws.Range("L21").Value = -0.14960882828101

$ws.Range("I22").Value = 0.07593777652216803
$ws.Range("J22").Value = 0.2704643720712596
$ws.Range("K22").Value = -0.1945265955490915
$ws.Range("L22").Value = -0.1752172027918623
